$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1: copy formatting from the adjacent header cell (E1, style "1")
# then overwrite its value with the new header text.
$ws.Cells.Item(1, 5).Copy($ws.Cells.Item(1, 6))
$ws.Cells.Item(1, 6).Value = "time_taken"

# New data column F2:F9 with plain (unstyled) text timestamps
$timestamps = @(
    "2021-10-05 13:40:04.487793",
    "2021-10-05 13:40:04.487803",
    "2021-10-05 13:40:04.487806",
    "2021-10-05 13:40:04.487809",
    "2021-10-05 13:40:04.487812",
    "2021-10-05 13:40:04.487815",
    "2021-10-05 13:40:04.487817",
    "2021-10-05 13:40:04.487820"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
